# The document ends with a paragraph containing a hyperlink to the MDN
# "Classes" reference page, followed by three empty paragraphs (the last of
# which carries the "_GoBack" bookmark). The edit removes the first of
# those three empty paragraphs, leaving two.

$d = $word.ActiveDocument

# Locate the paragraph that holds the MDN "Classes" hyperlink text.
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*developer.mozilla.org*Classes*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the anchor paragraph containing the MDN Classes hyperlink."
}

# The paragraph immediately following the anchor is the first of the
# trailing empty paragraphs; delete it (including its paragraph mark).
$emptyParaIndex = $anchorIndex + 1
$emptyPara = $d.Paragraphs.Item($emptyParaIndex)
$emptyPara.Range.Delete()

Write-Output "Deleted empty paragraph at index $emptyParaIndex (anchor was $anchorIndex)."
